# FIX: api uri modified
# The placeholder "API URI" style text in the first body paragraph is
# replaced with a new value.

$d = $word.ActiveDocument

$old = "회의_sdgfnlasdbviubasdvyiugsnuyuxasdfljbgksfdhjgkyfsdvkuziluasdfibsdhljgkvz"
$new = "sdnjfq3ruwhskdnfhauksdhauks"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
